$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.563.57"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "1.850.26"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.031"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.027"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4388"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07417"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "1.867.82"
$ws.Range("E12").Value = "  -8.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.516"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.692"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009050"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.027"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "27.597.00"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.269"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.915"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.976"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.284"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09053"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7618"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.527"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.881"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.029"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.150"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01977"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05308"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5160"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1677"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.790"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.516"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.713"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4655"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06399"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.849"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
